# "actualizadas diapos e informe"
# Slide 3 ("CSS"), content placeholder: merge the three runs that make up
# "Sitios web más fáciles de actualizar y modificar." into a single run,
# keeping the first run's formatting (lang="es-ES_tradnl" dirty="0").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Locate the paragraph that currently reads
# "Sitios " + "web más " + "fáciles de actualizar y modificar."
$para = $null
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $candidate = $tr.Paragraphs($i)
    if ($candidate.Text -like "Sitios*actualizar y modificar.*") {
        $para = $candidate
        break
    }
}

# Keep the leading "Sitios " (first run, carries the rPr we want to preserve)
# and fold the remainder of the paragraph into it as a single run.
$firstRun = $para.Characters(1, 7)
$remainder = $para.Characters(8, $para.Length - 7)
$remainder.Text = ""
$firstRun.InsertAfter("web más fáciles de actualizar y modificar.") | Out-Null
